$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.952.77"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "1.638.59"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "'215.54"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").Value = "1.866.12"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D14").Value = "1.639.17"
$ws.Range("E14").Value = "  -1.74%  "

$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "'63.00"
$ws.Range("E17").Value = "  -1.05%  "

$ws.Range("D18").Value = "26.067.59"
$ws.Range("E18").Value = "  -0.94%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("D22").Value = "'9.93"
$ws.Range("E22").Value = "  -1.64%  "

$ws.Range("D23").Value = "'6.29"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "'143.41"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  +3.76%  "

$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("D29").Value = "'15.59"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("D34").Value = "'1.54"
$ws.Range("E34").Value = "  -4.92%  "

$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").Value = "1.136.36"

$ws.Range("E38").Value = "  -2.08%  "

$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D42").Value = "'5.49"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("D43").Value = "'99.22"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("D44").Value = "'0.798"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").Value = "1.775.82"
$ws.Range("E45").Value = "  -0.53%  "

$ws.Range("E46").Value = "  +2.49%  "

$ws.Range("D47").Value = "'56.69"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("E48").Value = "  +2.58%  "

$ws.Range("D49").Value = "'1.47"
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("D50").Value = "'7.67"
$ws.Range("E50").Value = "  -0.26%  "
